$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1826.5588
$ws.Cells.Item(15, 9).Value = 1826.5588
$ws.Cells.Item(15, 11).Value = 5479.6764
$ws.Cells.Item(15, 13).Value = -5310.6764
$ws.Cells.Item(110, 8).Value = 702000000
$ws.Cells.Item(110, 10).Value = 702000000
$ws.Cells.Item(110, 12).Value = 702000000
$ws.Cells.Item(110, 14).Value = -702008180
$ws.Cells.Item(129, 8).Value = 111112650
$ws.Cells.Item(129, 9).Value = 924.5
$ws.Cells.Item(129, 10).Value = 200002020
$ws.Cells.Item(129, 11).Value = 2773.5
$ws.Cells.Item(129, 12).Value = 600006060
$ws.Cells.Item(129, 13).Value = 2226.5
$ws.Cells.Item(129, 14).Value = -600016060
$ws.Cells.Item(132, 8).Value = 2400.9207
$ws.Cells.Item(132, 9).Value = 2555.451
$ws.Cells.Item(132, 10).Value = 1744.1666
$ws.Cells.Item(132, 11).Value = 7666.353
$ws.Cells.Item(132, 12).Value = 5232.4998
$ws.Cells.Item(132, 13).Value = -5136.353
$ws.Cells.Item(132, 14).Value = -10292.4998
$ws.Cells.Item(141, 8).Value = 5422.087
$ws.Cells.Item(141, 9).Value = 5118.067
$ws.Cells.Item(141, 11).Value = 15354.201
$ws.Cells.Item(141, 13).Value = -10174.201

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4181.1665
$ws.Cells.Item(61, 9).Value = 4181.1665
$ws.Cells.Item(61, 11).Value = 4181.1665
$ws.Cells.Item(61, 13).Value = -3969.1665
$ws.Cells.Item(74, 8).Value = 2694.7058
$ws.Cells.Item(74, 9).Value = 2193
$ws.Cells.Item(74, 10).Value = 3898.8
$ws.Cells.Item(74, 11).Value = 2193
$ws.Cells.Item(74, 12).Value = 3898.8
$ws.Cells.Item(74, 13).Value = -1319
$ws.Cells.Item(74, 14).Value = -5646.8
$ws.Cells.Item(77, 8).Value = 2694.7058
$ws.Cells.Item(77, 9).Value = 2193
$ws.Cells.Item(77, 10).Value = 3898.8
$ws.Cells.Item(77, 11).Value = 10965
$ws.Cells.Item(77, 12).Value = 19494
$ws.Cells.Item(77, 13).Value = -6597
$ws.Cells.Item(77, 14).Value = -28230
$ws.Cells.Item(102, 8).Value = 3453.8667
$ws.Cells.Item(102, 9).Value = 2960.9167
$ws.Cells.Item(102, 11).Value = 2960.9167
$ws.Cells.Item(102, 13).Value = -1338.9167
$ws.Cells.Item(132, 8).Value = 2809.9033
$ws.Cells.Item(132, 9).Value = 2420.5925
$ws.Cells.Item(132, 11).Value = 7261.7775
$ws.Cells.Item(132, 13).Value = -4731.7775
$ws.Cells.Item(136, 8).Value = 4181.1665
$ws.Cells.Item(136, 9).Value = 4181.1665
$ws.Cells.Item(136, 11).Value = 12543.4995
$ws.Cells.Item(136, 13).Value = -9993.499500000002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2492.3438
$ws.Cells.Item(94, 9).Value = 2043
$ws.Cells.Item(94, 11).Value = 2043
$ws.Cells.Item(94, 13).Value = -1592
$ws.Cells.Item(99, 8).Value = 41401.89
$ws.Cells.Item(99, 10).Value = 7750
$ws.Cells.Item(99, 12).Value = 7750
$ws.Cells.Item(99, 14).Value = -10746
$ws.Cells.Item(134, 8).Value = 7243.8076
$ws.Cells.Item(134, 9).Value = 7333.56
$ws.Cells.Item(134, 10).Value = 5000
$ws.Cells.Item(134, 11).Value = 22000.68
$ws.Cells.Item(134, 12).Value = 15000
$ws.Cells.Item(134, 13).Value = -19465.68
$ws.Cells.Item(134, 14).Value = -20070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2292.2942
$ws.Cells.Item(16, 9).Value = 2347.9167
$ws.Cells.Item(16, 11).Value = 2347.9167
$ws.Cells.Item(16, 13).Value = -2060.9167
$ws.Cells.Item(31, 8).Value = 8373.4
$ws.Cells.Item(31, 9).Value = 7746.8
$ws.Cells.Item(31, 10).Value = 9000
$ws.Cells.Item(31, 11).Value = 7746.8
$ws.Cells.Item(31, 12).Value = 9000
$ws.Cells.Item(31, 13).Value = -7451.8
$ws.Cells.Item(31, 14).Value = -9590
$ws.Cells.Item(34, 8).Value = 8373.4
$ws.Cells.Item(34, 9).Value = 7746.8
$ws.Cells.Item(34, 10).Value = 9000
$ws.Cells.Item(34, 11).Value = 7746.8
$ws.Cells.Item(34, 12).Value = 9000
$ws.Cells.Item(34, 13).Value = -7544.8
$ws.Cells.Item(34, 14).Value = -9404
$ws.Cells.Item(62, 8).Value = 20666.666
$ws.Cells.Item(62, 9).Value = 25000
$ws.Cells.Item(62, 11).Value = 25000
$ws.Cells.Item(62, 13).Value = -24376
$ws.Cells.Item(65, 8).Value = 20666.666
$ws.Cells.Item(65, 9).Value = 25000
$ws.Cells.Item(65, 11).Value = 125000
$ws.Cells.Item(65, 13).Value = -121880
$ws.Cells.Item(113, 8).Value = 2292.2942
$ws.Cells.Item(113, 9).Value = 2347.9167
$ws.Cells.Item(113, 11).Value = 2347.9167
$ws.Cells.Item(113, 13).Value = -177.9167000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 172200
$ws.Cells.Item(80, 9).Value = 6000
$ws.Cells.Item(80, 11).Value = 18000
$ws.Cells.Item(80, 13).Value = -17064
$ws.Cells.Item(83, 8).Value = 172200
$ws.Cells.Item(83, 9).Value = 6000
$ws.Cells.Item(83, 11).Value = 54000
$ws.Cells.Item(83, 13).Value = -49320
$ws.Cells.Item(92, 8).Value = 469.25
$ws.Cells.Item(92, 9).Value = 475.66666
$ws.Cells.Item(92, 11).Value = 1426.99998
$ws.Cells.Item(92, 13).Value = -178.9999800000001
$ws.Cells.Item(129, 8).Value = 25642516
$ws.Cells.Item(129, 10).Value = 47621348
$ws.Cells.Item(129, 12).Value = 142864044
$ws.Cells.Item(129, 14).Value = -142874044
$ws.Cells.Item(140, 8).Value = 8097.593
$ws.Cells.Item(140, 9).Value = 8097.593
$ws.Cells.Item(140, 11).Value = 24292.779
$ws.Cells.Item(140, 13).Value = -19112.779
$ws.Cells.Item(141, 8).Value = 2267.8333
$ws.Cells.Item(141, 9).Value = 1821.4
$ws.Cells.Item(141, 10).Value = 4500
$ws.Cells.Item(141, 11).Value = 5464.200000000001
$ws.Cells.Item(141, 12).Value = 13500
$ws.Cells.Item(141, 13).Value = -284.2000000000007
$ws.Cells.Item(141, 14).Value = -23860

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2500.4
$ws.Cells.Item(80, 9).Value = 2477.4285
$ws.Cells.Item(80, 10).Value = 2520.5
$ws.Cells.Item(80, 11).Value = 2477.4285
$ws.Cells.Item(80, 12).Value = 2520.5
$ws.Cells.Item(80, 13).Value = -1479.4285
$ws.Cells.Item(80, 14).Value = -4516.5
$ws.Cells.Item(83, 8).Value = 2500.4
$ws.Cells.Item(83, 9).Value = 2477.4285
$ws.Cells.Item(83, 10).Value = 2520.5
$ws.Cells.Item(83, 11).Value = 12387.1425
$ws.Cells.Item(83, 12).Value = 12602.5
$ws.Cells.Item(83, 13).Value = -7395.1425
$ws.Cells.Item(83, 14).Value = -22586.5
$ws.Cells.Item(97, 8).Value = 7056.275
$ws.Cells.Item(97, 9).Value = 9057.679
$ws.Cells.Item(97, 10).Value = 2386.3333
$ws.Cells.Item(97, 11).Value = 9057.679
$ws.Cells.Item(97, 12).Value = 2386.3333
$ws.Cells.Item(97, 13).Value = -8561.679
$ws.Cells.Item(97, 14).Value = -3378.3333
$ws.Cells.Item(122, 8).Value = 16691.957
$ws.Cells.Item(122, 9).Value = 14762.368
$ws.Cells.Item(122, 11).Value = 44287.104
$ws.Cells.Item(122, 13).Value = -41837.104

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 271.18182
$ws.Cells.Item(16, 9).Value = 278.3
$ws.Cells.Item(16, 10).Value = 200
$ws.Cells.Item(16, 11).Value = 278.3
$ws.Cells.Item(16, 12).Value = 200
$ws.Cells.Item(16, 13).Value = -108.3
$ws.Cells.Item(16, 14).Value = -540
$ws.Cells.Item(55, 8).Value = 962.6799999999999
$ws.Cells.Item(55, 9).Value = 882.8946999999999
$ws.Cells.Item(55, 11).Value = 882.8946999999999
$ws.Cells.Item(55, 13).Value = -709.8946999999999
$ws.Cells.Item(68, 8).Value = 5452.0625
$ws.Cells.Item(68, 9).Value = 3822
$ws.Cells.Item(68, 10).Value = 5684.9287
$ws.Cells.Item(68, 11).Value = 3822
$ws.Cells.Item(68, 12).Value = 5684.9287
$ws.Cells.Item(68, 13).Value = -3073
$ws.Cells.Item(68, 14).Value = -7182.9287
$ws.Cells.Item(71, 8).Value = 5452.0625
$ws.Cells.Item(71, 9).Value = 3822
$ws.Cells.Item(71, 10).Value = 5684.9287
$ws.Cells.Item(71, 11).Value = 19110
$ws.Cells.Item(71, 12).Value = 28424.6435
$ws.Cells.Item(71, 13).Value = -15366
$ws.Cells.Item(71, 14).Value = -35912.64350000001
$ws.Cells.Item(132, 8).Value = 1501665.8
$ws.Cells.Item(132, 9).Value = 1876082.9
$ws.Cells.Item(132, 10).Value = 3997.5
$ws.Cells.Item(132, 11).Value = 5628248.699999999
$ws.Cells.Item(132, 12).Value = 11992.5
$ws.Cells.Item(132, 13).Value = -5625718.699999999
$ws.Cells.Item(132, 14).Value = -17052.5
$ws.Cells.Item(136, 8).Value = 7274.4443
$ws.Cells.Item(136, 9).Value = 4199.8
$ws.Cells.Item(136, 10).Value = 8457
$ws.Cells.Item(136, 11).Value = 12599.4
$ws.Cells.Item(136, 12).Value = 25371
$ws.Cells.Item(136, 13).Value = -10049.4
$ws.Cells.Item(136, 14).Value = -30471

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 54645
$ws.Cells.Item(122, 10).Value = 96665.664
$ws.Cells.Item(122, 12).Value = 289996.992
$ws.Cells.Item(122, 14).Value = -294896.992
$ws.Cells.Item(132, 8).Value = 11860.23
$ws.Cells.Item(132, 9).Value = 13114.552
$ws.Cells.Item(132, 10).Value = 8222.700000000001
$ws.Cells.Item(132, 11).Value = 39343.656
$ws.Cells.Item(132, 12).Value = 24668.1
$ws.Cells.Item(132, 13).Value = -36813.656
$ws.Cells.Item(132, 14).Value = -29728.1
$ws.Cells.Item(136, 8).Value = 2147.5
$ws.Cells.Item(136, 9).Value = 1486.2222
$ws.Cells.Item(136, 10).Value = 8099
$ws.Cells.Item(136, 11).Value = 4458.6666
$ws.Cells.Item(136, 12).Value = 24297
$ws.Cells.Item(136, 13).Value = -1908.6666
$ws.Cells.Item(136, 14).Value = -29397
